$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The edit effectively *moves* the "Meta description" paragraph (currently
# paragraph #2, right after the H1 title) down to the very end of the
# document, turning it into two new paragraphs that replace the old
# "Prompt for DALLE..." paragraph's text:
#   1) a bold paragraph repeating the title text
#   2) an italic paragraph with the (former) meta-description body text
# The old DALLE paragraph's own run/paragraph structure is reused (via
# Find/Replace) so its formatting (leading empty run + italic run) carries
# over untouched.
# ---------------------------------------------------------------------------

# Step 1: clone the "Meta description" paragraph's exact run structure
# (leading empty run + bold "Meta description" run + plain trailing run)
# to a point right before the last paragraph (the DALLE prompt paragraph).
$n = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($n)
$insertPoint = $d.Range($pLast.Range.Start, $pLast.Range.Start)

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Copy()
$insertPoint.Paste()

# Step 2: turn the freshly-pasted clone's text ("Meta description: Experience
# ...") into the headline text, reusing its bold run's formatting.
$d.Content.Find.Execute("Meta description: Experience the innovative gameplay of Cool Jewels. Read our review and play for free.", $true, $false, $false, $false, $false, $true, 1, $false, "Play Cool Jewels Free - Review of Innovative Slot Design", 2)

# Step 3: turn the original last paragraph's DALLE-prompt text into the
# meta-description body text, reusing its own (italic) run formatting.
$d.Content.Find.Execute("Prompt for DALLE: Create a feature image for Cool Jewels that captures the game's fun and engaging gameplay. The image should be in a cartoon style and prominently feature a happy Maya warrior with glasses. The warrior should be surrounded by colorful gems and jewels, with lightning bolts striking some of them. Make sure to incorporate the game's name `"Cool Jewels`" in the design as well.", $true, $false, $false, $false, $false, $true, 1, $false, "Experience the innovative gameplay of Cool Jewels. Read our review and play for free.", 2)

# Step 4: remove the original "Meta description" paragraph near the top of
# the document (right after the H1 title paragraph).
$d.Paragraphs.Item(2).Range.Delete()
